$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 31 (this shifts old rows 31-35 down to 32-36,
# preserving all their cell formatting).
$ws.Rows.Item(31).EntireRow.Insert()

# Populate the newly inserted row 31 with this week's data (same structure/unit/
# origin/classification as the rest of the Jengibre series, new date + figures).
$ws.Cells.Item(31,1).Value = 11
$ws.Cells.Item(31,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(31,3).Value = "Bíobío"
$ws.Cells.Item(31,4).Value = 44769
$ws.Cells.Item(31,5).Value = 8
$ws.Cells.Item(31,6).Value = 100114007
$ws.Cells.Item(31,7).Value = "Jengibre"
$ws.Cells.Item(31,8).Value = "Sin especificar"
$ws.Cells.Item(31,9).Value = "Primera"
$ws.Cells.Item(31,10).Value = 50
$ws.Cells.Item(31,11).Value = 14000
$ws.Cells.Item(31,12).Value = 15000
$ws.Cells.Item(31,13).Value = 14600
$ws.Cells.Item(31,14).Value = "$/caja 13 kilos"
$ws.Cells.Item(31,15).Value = "Perú"
$ws.Cells.Item(31,16).Value = 1123
$ws.Cells.Item(31,17).Value = 13
$ws.Cells.Item(31,18).Value = "Hortaliza"
